# "Pie chart output with matplotlib"
# The underlying data was regenerated, introducing two new "function"
# categories ("redundancy" and "reproduction") that replace the previous
# value in a handful of rows of column D ("function").

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$updates = @{
    "D3"  = "redundancy"
    "D4"  = "reproduction"
    "D9"  = "redundancy"
    "D11" = "reproduction"
    "D14" = "redundancy"
    "D15" = "reproduction"
    "D17" = "redundancy"
    "D29" = "reproduction"
    "D40" = "reproduction"
    "D41" = "redundancy"
    "D42" = "redundancy"
    "D45" = "reproduction"
}

foreach ($addr in $updates.Keys) {
    $ws.Range($addr).Value = $updates[$addr]
}

# Restore the view/selection state recorded in the workbook: the author had
# scrolled the sheet down and left the active cell on F39.
$win = $excel.ActiveWindow
$ws.Range("F39").Select()
$win.ScrollRow = 29
$win.ScrollColumn = 1

$wb.Save()
